# Consolidate text runs that are split across multiple <a:r> elements but
# whose concatenated text did not actually need to change, e.g.
#   "A" + " " + "slide"              -> "A slide"
#   "Just" + " " + "an" + ... "side" -> "Just an image on this side"
#
# Simply assigning .TextFrame.TextRange.Text to the same (already equal)
# string is treated as a no-op by the writer, and assigning it to a
# genuinely different string creates a fresh run with default formatting
# (e.g. adds lang="en-US"), which would not match the original run
# properties. Re-assigning the text through .Characters(1, Length),
# however, rewrites the whole range as a single run while preserving the
# existing (empty) run properties, so use that to force the consolidation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Characters(1, $title.Length).Text = "A slide"

$caption = $s.Shapes.Item(4).TextFrame.TextRange
$caption.Characters(1, $caption.Length).Text = "Just an image on this side"
